$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.605.31'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '3.121.87'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.34%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.120.53'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.36'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = '3.632.99'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.120'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").Value = '63.675.08'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '3.116.31'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.707'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("E25").Value = '  -3.34%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -3.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.58%  '
$ws.Range("E30").Value = '  -3.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.111'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  -6.50%  '
$ws.Range("E39").Value = '  -7.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '438.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0393'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.91%  '
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("D44").Value = '2.865.53'
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.259'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.12%  '
$ws.Range("E46").Value = '  -4.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.114'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.07%  '
